$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the workbook's new sheet name
$ws.Name = "ConfigurationInputsTrout"

# Add the new "#ProductionPeriod" section with ProdStartDay / ProdEndDay rows
$ws.Range("A27").Value = "#ProductionPeriod"

$ws.Range("A28").Value = "ProdStartDay"
$ws.Range("B28").Value = 105
$ws.Range("C28").Value = "JulianDay"

$ws.Range("A29").Value = "ProdEndDay"
$ws.Range("B29").Value = 330
$ws.Range("C29").Value = "JulianDay"

# Match the selection left by the author after entering the new data
$ws.Range("B30").Select()
